$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.389221
$ws.Range("H2").Value = 4.167663
$ws.Range("I2").Value = 0.2910270461264192
$ws.Range("J2").Value = 0.2910270461264192
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 202.4105342065714
$ws.Range("R2").Value = 1821.694807859142
$ws.Range("S2").Value = 0.0834065397795499
$ws.Range("T2").Value = 0.0834065397795499
$ws.Range("G3").Value = 1.389221
$ws.Range("H3").Value = 4.167663
$ws.Range("I3").Value = 0.2910270461264192
$ws.Range("J3").Value = 0.2910270461264192
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 234.5000917382894
$ws.Range("R3").Value = 2110.500825644604
$ws.Range("S3").Value = 0.09662956182861922
$ws.Range("T3").Value = 0.09662956182861922
$ws.Range("G4").Value = 1.389221
$ws.Range("H4").Value = 4.167663
$ws.Range("I4").Value = 0.2910270461264192
$ws.Range("J4").Value = 0.2910270461264192
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 177.9955447121814
$ws.Range("R4").Value = 1601.959902409632
$ws.Range("S4").Value = 0.07334594782239957
$ws.Range("T4").Value = 0.07334594782239957
$ws.Range("G5").Value = 1.389221
$ws.Range("H5").Value = 4.167663
$ws.Range("I5").Value = 0.2910270461264192
$ws.Range("J5").Value = 0.2910270461264192
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 91.356671929459
$ws.Range("R5").Value = 822.210047365131
$ws.Range("S5").Value = 0.03764499669585049
$ws.Range("T5").Value = 0.03764499669585049
$ws.Range("I6").Value = 0.461328155686921
$ws.Range("J6").Value = 0.4613281556869209
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 320.8556719383383
$ws.Range("R6").Value = 2887.701047445045
$ws.Range("S6").Value = 0.1322137776569852
$ws.Range("T6").Value = 0.1322137776569852
$ws.Range("I7").Value = 0.461328155686921
$ws.Range("J7").Value = 0.4613281556869209
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.1531745524567775
$ws.Range("T7").Value = 0.1531745524567775
$ws.Range("I8").Value = 0.461328155686921
$ws.Range("J8").Value = 0.4613281556869209
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 282.1536948386916
$ws.Range("R8").Value = 2539.383253548224
$ws.Range("S8").Value = 0.1162660010001905
$ws.Range("T8").Value = 0.1162660010001905
$ws.Range("I9").Value = 0.461328155686921
$ws.Range("J9").Value = 0.4613281556869209
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 144.8161108456047
$ws.Range("R9").Value = 1303.344997610442
$ws.Range("S9").Value = 0.05967382457296778
$ws.Range("T9").Value = 0.05967382457296778
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1506176666666667
$ws.Range("H10").Value = 0.451853
$ws.Range("I10").Value = 0.03155280162368235
$ws.Range("J10").Value = 0.03155280162368235
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 21.94510619328911
$ws.Range("R10").Value = 197.505955739602
$ws.Range("S10").Value = 0.009042836529491219
$ws.Range("T10").Value = 0.009042836529491219
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1506176666666667
$ws.Range("H11").Value = 0.451853
$ws.Range("I11").Value = 0.03155280162368235
$ws.Range("J11").Value = 0.03155280162368235
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 25.42421734968045
$ws.Range("R11").Value = 228.817956147124
$ws.Range("S11").Value = 0.01047646064495788
$ws.Range("T11").Value = 0.01047646064495788
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1506176666666667
$ws.Range("H12").Value = 0.451853
$ws.Range("I12").Value = 0.03155280162368235
$ws.Range("J12").Value = 0.03155280162368235
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 19.29806245486578
$ws.Range("R12").Value = 173.682562093792
$ws.Range("S12").Value = 0.007952079273538843
$ws.Range("T12").Value = 0.007952079273538843
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1506176666666667
$ws.Range("H13").Value = 0.451853
$ws.Range("I13").Value = 0.03155280162368235
$ws.Range("J13").Value = 0.03155280162368235
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 9.904780276462333
$ws.Range("R13").Value = 89.14302248816099
$ws.Range("S13").Value = 0.00408142517569442
$ws.Range("T13").Value = 0.00408142517569442
$ws.Range("G14").Value = 1.031517666666667
$ws.Range("H14").Value = 3.094553
$ws.Range("I14").Value = 0.2160919965629775
$ws.Range("J14").Value = 0.2160919965629775
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 150.2928921701558
$ws.Range("R14").Value = 1352.636029531402
$ws.Range("S14").Value = 0.06193062104455793
$ws.Range("T14").Value = 0.06193062104455793
$ws.Range("G15").Value = 1.031517666666667
$ws.Range("H15").Value = 3.094553
$ws.Range("I15").Value = 0.2160919965629775
$ws.Range("J15").Value = 0.2160919965629775
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 174.1198754287471
$ws.Range("R15").Value = 1567.078878858724
$ws.Range("S15").Value = 0.07174891550622953
$ws.Range("T15").Value = 0.07174891550622953
$ws.Range("G16").Value = 1.031517666666667
$ws.Range("H16").Value = 3.094553
$ws.Range("I16").Value = 0.2160919965629775
$ws.Range("J16").Value = 0.2160919965629775
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 132.1643921007324
$ws.Range("R16").Value = 1189.479528906592
$ws.Range("S16").Value = 0.05446047889948157
$ws.Range("T16").Value = 0.05446047889948157
$ws.Range("G17").Value = 1.031517666666667
$ws.Range("H17").Value = 3.094553
$ws.Range("I17").Value = 0.2160919965629775
$ws.Range("J17").Value = 0.2160919965629775
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 67.83371476756233
$ws.Range("R17").Value = 610.503432908061
$ws.Range("S17").Value = 0.02795198111270854
$ws.Range("T17").Value = 0.02795198111270854
